$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.985.15"
$ws.Range("E2").Value = "  +11.10%  "

$ws.Range("D3").Value = "1.808.12"
$ws.Range("E3").Value = "  +7.63%  "

$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").Value = "'227.32"
$ws.Range("E5").Value = "  +3.38%  "

$ws.Range("D6").Value = "'0.542"
$ws.Range("E6").Value = "  +3.23%  "

$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("D8").Value = "'31.30"
$ws.Range("E8").Value = "  +4.10%  "

$ws.Range("D9").Value = "'47.09"
$ws.Range("E9").Value = "  +6.54%  "

$ws.Range("D10").Value = "'0.280"
$ws.Range("E10").Value = "  +6.00%  "

$ws.Range("D11").Value = "'0.0665"
$ws.Range("E11").Value = "  +6.12%  "

$ws.Range("D12").Value = "'0.0926"
$ws.Range("E12").Value = "  +2.06%  "

$ws.Range("D13").Value = "2.065.86"
$ws.Range("E13").Value = "  +7.52%  "

$ws.Range("D14").Value = "1.812.96"
$ws.Range("E14").Value = "  +6.71%  "

$ws.Range("D15").Value = "'0.636"
$ws.Range("E15").Value = "  +2.49%  "

$ws.Range("D16").Value = "33.918.75"
$ws.Range("E16").Value = "  +10.90%  "

$ws.Range("D17").Value = "'10.16"
$ws.Range("E17").Value = "  -3.29%  "

$ws.Range("D18").Value = "'4.25"
$ws.Range("E18").Value = "  +6.93%  "

$ws.Range("D19").Value = "'69.10"
$ws.Range("E19").Value = "  +4.30%  "

$ws.Range("D20").Value = "'256.43"
$ws.Range("E20").Value = "  +4.85%  "

$ws.Range("D21").Value = "0.0₃0745"
$ws.Range("E21").Value = "  +4.03%  "

$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").Value = "'10.47"
$ws.Range("E23").Value = "  +2.84%  "

$ws.Range("D24").Value = "'4.32"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("D26").Value = "'157.82"
$ws.Range("E26").Value = "  +0.02%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'16.53"
$ws.Range("E27").Value = "  +4.25%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'7.07"
$ws.Range("E28").Value = "  +5.44%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.115"
$ws.Range("E29").Value = "  +3.10%  "

$ws.Range("B30").Value = "MinaProtocolToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D30").Value = "'2.20"
$ws.Range("E30").Value = "  +435.43%  "

$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("D32").Value = "'3.85"
$ws.Range("E32").Value = "  +10.73%  "

$ws.Range("D33").Value = "'0.0511"
$ws.Range("E33").Value = "  +2.77%  "

$ws.Range("E34").Value = "  +4.73%  "

$ws.Range("D35").Value = "'3.50"
$ws.Range("E35").Value = "  +6.43%  "

$ws.Range("D36").Value = "1.536.33"
$ws.Range("E36").Value = "  +1.83%  "

$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  +2.43%  "

$ws.Range("D38").Value = "'1.07"
$ws.Range("E38").Value = "  +3.75%  "

$ws.Range("D39").Value = "'84.18"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "'0.0187"
$ws.Range("E40").Value = "  +4.68%  "

$ws.Range("D41").Value = "'0.617"
$ws.Range("E41").Value = "  +5.04%  "

$ws.Range("E42").Value = "  +3.30%  "

$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("D44").Value = "'0.907"
$ws.Range("E44").Value = "  +8.27%  "

$ws.Range("D45").Value = "'2.12"
$ws.Range("E45").Value = "  +6.61%  "

$ws.Range("E46").Value = "  +4.23%  "

$ws.Range("E47").Value = "  +4.28%  "

$ws.Range("D48").Value = "1.965.85"
$ws.Range("E48").Value = "  +8.08%  "

$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("D50").Value = "'5.67"
$ws.Range("E50").Value = "  +2.99%  "

$ws.Range("D51").Value = "'52.54"
$ws.Range("E51").Value = "  +1.98%  "

